# Apply the updates described by the commit diff.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # sheet "23"
$ws2 = $wb.Worksheets.Item(2)   # sheet "34"

# --- Sheet "23" value edits ---------------------------------------------
# Row 3
$ws1.Range("B3").Value = 518
$ws1.Range("C3").Value = 5.5
$ws1.Range("D3").Value = 543
$ws1.Range("E3").Value = 52

# Row 4
$ws1.Range("B4").Value = 7
$ws1.Range("C4").Value = 314.67
$ws1.Range("D4").Value = 32
$ws1.Range("E4").Value = 361.17

# Row 5
$ws1.Range("B5").Value = 519
$ws1.Range("C5").Value = 307.67
$ws1.Range("D5").Value = 544
$ws1.Range("E5").Value = 354.17

# Row 18
$ws1.Range("B18").Value = 69.355199999999996
$ws1.Range("C18").Value = 128.75749999999999
$ws1.Range("D18").Value = 203.3552
$ws1.Range("E18").Value = 162.75749999999999

# --- Sheet "34" value edits ----------------------------------------------
# Row 18
$ws2.Range("B18").Value = 38.799700000000001
$ws2.Range("C18").Value = 128.75749999999999
$ws2.Range("D18").Value = 172.7997
$ws2.Range("E18").Value = 162.75749999999999

# --- Selection / active-cell updates ------------------------------------
# Set sheet "34" selection first, then re-activate sheet "23" so it stays
# the tab shown when the workbook is opened (matches the original state).
$ws2.Activate()
$ws2.Range("E28").Select()

$ws1.Activate()
$ws1.Range("D27").Select()
